{"js": "// Update each piece of text in the document (the date line and the\n// division-problem table cells) to the new values, preserving all\n// existing run formatting by doing an in-place text replace.\nconst replacements = [\n  [\"2024-12-28 Saturday\", \"2024-12-29 Sunday\"],\n  [\"458\u00f72=229, 0\", \"691\u00f76=115, 1\"],\n  [\"511\u00f75=102, 1\", \"412\u00f78=51, 4\"],\n  [\"150\u00f73=50, 0\", \"368\u00f75=73, 3\"],\n  [\"941\u00f73=313, 2\", \"764\u00f76=127, 2\"],\n  [\"920\u00f75=184, 0\", \"426\u00f76=71, 0\"],\n  [\"217\u00f76=36, 1\", \"152\u00f72=76, 0\"],\n  [\"967\u00f74=241, 3\", \"296\u00f79=32, 8\"],\n  [\"384\u00f77=54, 6\", \"727\u00f74=181, 3\"],\n  [\"597\u00f76=99, 3\", \"711\u00f77=101, 4\"],\n  [\"306\u00f76=51, 0\", \"884\u00f79=98, 2\"],\n  [\"562\u00f75=112, 2\", \"774\u00f75=154, 4\"],\n  [\"494\u00f77=70, 4\", \"864\u00f75=172, 4\"],\n  [\"459\u00f76=76, 3\", \"740\u00f79=82, 2\"],\n  [\"265\u00f79=29, 4\", \"798\u00f72=399, 0\"],\n  [\"794\u00f77=113, 3\", \"697\u00f77=99, 4\"],\n  [\"363\u00f75=72, 3\", \"329\u00f72=164, 1\"],\n  [\"657\u00f73=219, 0\", \"958\u00f79=106, 4\"],\n  [\"985\u00f73=328, 1\", \"944\u00f75=188, 4\"],\n  [\"766\u00f73=255, 1\", \"267\u00f78=33, 3\"],\n  [\"407\u00f79=45, 2\", \"865\u00f79=96, 1\"],\n  [\"206\u00f78=25, 6\", \"174\u00f73=58, 0\"],\n  [\"151\u00f72=75, 1\", \"598\u00f73=199, 1\"],\n  [\"185\u00f72=92, 1\", \"561\u00f78=70, 1\"],\n  [\"369\u00f77=52, 5\", \"483\u00f77=69, 0\"],\n  [\"210\u00f73=70, 0\", \"619\u00f76=103, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and each division-problem table cell to the new\n# values using Find/Replace, which preserves existing run formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-28 Saturday\", \"2024-12-29 Sunday\"),\n    @(\"458\u00f72=229, 0\", \"691\u00f76=115, 1\"),\n    @(\"511\u00f75=102, 1\", \"412\u00f78=51, 4\"),\n    @(\"150\u00f73=50, 0\", \"368\u00f75=73, 3\"),\n    @(\"941\u00f73=313, 2\", \"764\u00f76=127, 2\"),\n    @(\"920\u00f75=184, 0\", \"426\u00f76=71, 0\"),\n    @(\"217\u00f76=36, 1\", \"152\u00f72=76, 0\"),\n    @(\"967\u00f74=241, 3\", \"296\u00f79=32, 8\"),\n    @(\"384\u00f77=54, 6\", \"727\u00f74=181, 3\"),\n    @(\"597\u00f76=99, 3\", \"711\u00f77=101, 4\"),\n    @(\"306\u00f76=51, 0\", \"884\u00f79=98, 2\"),\n    @(\"562\u00f75=112, 2\", \"774\u00f75=154, 4\"),\n    @(\"494\u00f77=70, 4\", \"864\u00f75=172, 4\"),\n    @(\"459\u00f76=76, 3\", \"740\u00f79=82, 2\"),\n    @(\"265\u00f79=29, 4\", \"798\u00f72=399, 0\"),\n    @(\"794\u00f77=113, 3\", \"697\u00f77=99, 4\"),\n    @(\"363\u00f75=72, 3\", \"329\u00f72=164, 1\"),\n    @(\"657\u00f73=219, 0\", \"958\u00f79=106, 4\"),\n    @(\"985\u00f73=328, 1\", \"944\u00f75=188, 4\"),\n    @(\"766\u00f73=255, 1\", \"267\u00f78=33, 3\"),\n    @(\"407\u00f79=45, 2\", \"865\u00f79=96, 1\"),\n    @(\"206\u00f78=25, 6\", \"174\u00f73=58, 0\"),\n    @(\"151\u00f72=75, 1\", \"598\u00f73=199, 1\"),\n    @(\"185\u00f72=92, 1\", \"561\u00f78=70, 1\"),\n    @(\"369\u00f77=52, 5\", \"483\u00f77=69, 0\"),\n    @(\"210\u00f73=70, 0\", \"619\u00f76=103, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
